$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column F (ASSISTS) rows 2-41 were stored as inline text ("0", "2", ...)
# but should be stored as real numbers with the same value. Re-assigning
# the numeric value forces Excel to persist them as numeric cells.
$assists = @{
    2  = 0
    3  = 0
    4  = 0
    5  = 0
    6  = 0
    7  = 2
    8  = 2
    9  = 2
    10 = 2
    11 = 2
    12 = 2
    13 = 2
    14 = 4
    15 = 4
    16 = 4
    17 = 4
    18 = 5
    19 = 5
    20 = 5
    21 = 5
    22 = 5
    23 = 5
    24 = 5
    25 = 6
    26 = 6
    27 = 6
    28 = 7
    29 = 7
    30 = 9
    31 = 10
    32 = 10
    33 = 10
    34 = 10
    35 = 10
    36 = 10
    37 = 11
    38 = 11
    39 = 11
    40 = 11
    41 = 11
}

foreach ($row in $assists.Keys) {
    $ws.Cells.Item($row, 6).Value = $assists[$row]
}

# Column H (CHAMPION) fixes: several rows had the wrong champion recorded
# and should all read "Galio".
$championRows = @(5, 11, 17, 21, 23, 28, 29, 30, 35, 37, 38, 41)
foreach ($row in $championRows) {
    $ws.Cells.Item($row, 8).Value = "Galio"
}
